$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1136.5454
$ws.Range("I100").Value = 1136.5454
$ws.Range("K100").Value = 1136.5454
$ws.Range("M100").Value = -595.5454
$ws.Range("H129").Value = 2761.9622
$ws.Range("J129").Value = 940.8222
$ws.Range("L129").Value = 2822.4666
$ws.Range("N129").Value = -12822.4666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 67663.336
$ws.Range("I2").Value = 1110.8462
$ws.Range("J2").Value = 500254.5
$ws.Range("K2").Value = 1110.8462
$ws.Range("L2").Value = 500254.5
$ws.Range("M2").Value = -997.8462
$ws.Range("N2").Value = -500480.5
$ws.Range("H45").Value = 56675.277
$ws.Range("I45").Value = 91528.27
$ws.Range("J45").Value = 1906.2858
$ws.Range("K45").Value = 91528.27
$ws.Range("L45").Value = 1906.2858
$ws.Range("M45").Value = -91151.27
$ws.Range("N45").Value = -2660.2858
$ws.Range("H74").Value = 666.9184
$ws.Range("I74").Value = 571.6591
$ws.Range("J74").Value = 1505.2
$ws.Range("K74").Value = 571.6591
$ws.Range("L74").Value = 1505.2
$ws.Range("M74").Value = 302.3409
$ws.Range("N74").Value = -3253.2
$ws.Range("H77").Value = 666.9184
$ws.Range("I77").Value = 571.6591
$ws.Range("J77").Value = 1505.2
$ws.Range("K77").Value = 2858.2955
$ws.Range("L77").Value = 7526
$ws.Range("M77").Value = 1509.7045
$ws.Range("N77").Value = -16262
$ws.Range("H116").Value = 67663.336
$ws.Range("I116").Value = 1110.8462
$ws.Range("J116").Value = 500254.5
$ws.Range("K116").Value = 1110.8462
$ws.Range("L116").Value = 500254.5
$ws.Range("M116").Value = 1183.1538
$ws.Range("N116").Value = -504842.5
$ws.Range("H122").Value = 1885.6666
$ws.Range("I122").Value = 1541.375
$ws.Range("K122").Value = 4624.125
$ws.Range("M122").Value = -2174.125
$ws.Range("H132").Value = 18122.459
$ws.Range("I132").Value = 21675.232
$ws.Range("J132").Value = 2896.2856
$ws.Range("K132").Value = 65025.696
$ws.Range("L132").Value = 8688.856800000001
$ws.Range("M132").Value = -62495.696
$ws.Range("N132").Value = -13748.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 67663.336
$ws.Range("I3").Value = 1110.8462
$ws.Range("J3").Value = 500254.5
$ws.Range("K3").Value = 1110.8462
$ws.Range("L3").Value = 500254.5
$ws.Range("M3").Value = -996.8462
$ws.Range("N3").Value = -500482.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35020.91
$ws.Range("I31").Value = 1276.5
$ws.Range("J31").Value = 63141.25
$ws.Range("K31").Value = 1276.5
$ws.Range("L31").Value = 63141.25
$ws.Range("M31").Value = -981.5
$ws.Range("N31").Value = -63731.25
$ws.Range("H34").Value = 35020.91
$ws.Range("I34").Value = 1276.5
$ws.Range("J34").Value = 63141.25
$ws.Range("K34").Value = 1276.5
$ws.Range("L34").Value = 63141.25
$ws.Range("M34").Value = -1074.5
$ws.Range("N34").Value = -63545.25
$ws.Range("H134").Value = 965.5333000000001
$ws.Range("I134").Value = 653.2083
$ws.Range("K134").Value = 1959.6249
$ws.Range("M134").Value = 575.3751

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 45.6
$ws.Range("I12").Value = 17
$ws.Range("J12").Value = 50.64706
$ws.Range("K12").Value = 51
$ws.Range("L12").Value = 151.94118
$ws.Range("M12").Value = 122
$ws.Range("N12").Value = -497.94118
$ws.Range("H23").Value = 430.5263
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 453.33334
$ws.Range("K23").Value = 60
$ws.Range("L23").Value = 1360.00002
$ws.Range("M23").Value = 175
$ws.Range("N23").Value = -1830.00002
$ws.Range("H24").Value = 891.6667
$ws.Range("I24").Value = 450
$ws.Range("J24").Value = 1112.5
$ws.Range("K24").Value = 1350
$ws.Range("L24").Value = 3337.5
$ws.Range("M24").Value = -1120
$ws.Range("N24").Value = -3797.5
$ws.Range("H34").Value = 770.25
$ws.Range("I34").Value = 90
$ws.Range("J34").Value = 867.4286
$ws.Range("K34").Value = 270
$ws.Range("L34").Value = 2602.2858
$ws.Range("M34").Value = -186
$ws.Range("N34").Value = -2770.2858
$ws.Range("H39").Value = 3249.75
$ws.Range("J39").Value = 3642.5715
$ws.Range("L39").Value = 10927.7145
$ws.Range("N39").Value = -11515.7145
$ws.Range("H58").Value = 2200.7273
$ws.Range("I58").Value = 868.3333
$ws.Range("K58").Value = 2604.9999
$ws.Range("M58").Value = -2476.9999
$ws.Range("H81").Value = 1786.5385
$ws.Range("I81").Value = 906.6667
$ws.Range("J81").Value = 2252.353
$ws.Range("K81").Value = 2720.0001
$ws.Range("L81").Value = 6757.059
$ws.Range("M81").Value = -1597.0001
$ws.Range("N81").Value = -9003.059000000001
$ws.Range("H84").Value = 1786.5385
$ws.Range("I84").Value = 906.6667
$ws.Range("J84").Value = 2252.353
$ws.Range("K84").Value = 8160.0003
$ws.Range("L84").Value = 20271.177
$ws.Range("M84").Value = -2544.0003
$ws.Range("N84").Value = -31503.177
$ws.Range("H87").Value = 7744.857
$ws.Range("I87").Value = 4771.3335
$ws.Range("J87").Value = 9975
$ws.Range("K87").Value = 14314.0005
$ws.Range("L87").Value = 29925
$ws.Range("M87").Value = -13066.0005
$ws.Range("N87").Value = -32421
$ws.Range("H90").Value = 7744.857
$ws.Range("I90").Value = 4771.3335
$ws.Range("J90").Value = 9975
$ws.Range("K90").Value = 42942.0015
$ws.Range("L90").Value = 89775
$ws.Range("M90").Value = -36702.0015
$ws.Range("N90").Value = -102255
$ws.Range("H131").Value = 603392.3
$ws.Range("I131").Value = 660.8333
$ws.Range("J131").Value = 703847.5600000001
$ws.Range("K131").Value = 1982.4999
$ws.Range("L131").Value = 2111542.68
$ws.Range("M131").Value = 3057.5001
$ws.Range("N131").Value = -2121622.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 37594.4
$ws.Range("I34").Value = 20000
$ws.Range("J34").Value = 41993
$ws.Range("K34").Value = 20000
$ws.Range("L34").Value = 41993
$ws.Range("M34").Value = -19732
$ws.Range("N34").Value = -42529
$ws.Range("H64").Value = 300000
$ws.Range("I64").Value = 300000
$ws.Range("K64").Value = 300000
$ws.Range("M64").Value = -299752
$ws.Range("H67").Value = 300000
$ws.Range("I67").Value = 300000
$ws.Range("K67").Value = 300000
$ws.Range("M67").Value = -299142
$ws.Range("H76").Value = 37594.4
$ws.Range("I76").Value = 20000
$ws.Range("J76").Value = 41993
$ws.Range("K76").Value = 20000
$ws.Range("L76").Value = 41993
$ws.Range("M76").Value = -19685
$ws.Range("N76").Value = -42623
$ws.Range("H79").Value = 37594.4
$ws.Range("I79").Value = 20000
$ws.Range("J79").Value = 41993
$ws.Range("K79").Value = 20000
$ws.Range("L79").Value = 41993
$ws.Range("M79").Value = -18908
$ws.Range("N79").Value = -44177
$ws.Range("H92").Value = 18060.4
$ws.Range("J92").Value = 18060.4
$ws.Range("L92").Value = 18060.4
$ws.Range("N92").Value = -21804.4
$ws.Range("H93").Value = 34050
$ws.Range("J93").Value = 34050
$ws.Range("L93").Value = 34050
$ws.Range("N93").Value = -37794
$ws.Range("H122").Value = 2988.7368
$ws.Range("I122").Value = 3200.4614
$ws.Range("J122").Value = 2530
$ws.Range("K122").Value = 9601.3842
$ws.Range("L122").Value = 7590
$ws.Range("M122").Value = -7151.3842
$ws.Range("N122").Value = -12490
$ws.Range("H126").Value = 2649.8
$ws.Range("I126").Value = 2302.75
$ws.Range("J126").Value = 4038
$ws.Range("K126").Value = 6908.25
$ws.Range("L126").Value = 12114
$ws.Range("M126").Value = -4438.25
$ws.Range("N126").Value = -17054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 15000
$ws.Range("J97").Value = 15000
$ws.Range("L97").Value = 15000
$ws.Range("N97").Value = -16982
$ws.Range("H122").Value = 3733.6
$ws.Range("I122").Value = 3688
$ws.Range("J122").Value = 3785.7144
$ws.Range("K122").Value = 11064
$ws.Range("L122").Value = 11357.1432
$ws.Range("M122").Value = -8614
$ws.Range("N122").Value = -16257.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 17403.375
$ws.Range("I132").Value = 16105.1
$ws.Range("J132").Value = 19567.166
$ws.Range("K132").Value = 48315.3
$ws.Range("L132").Value = 58701.49800000001
$ws.Range("M132").Value = -45785.3
$ws.Range("N132").Value = -63761.49800000001
